$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the B2:D9 range to 0 by default
$ws.Range("B2:D9").Value = 0

# Override the two cells that take on new non-zero values
$ws.Range("B4").Value = -0.8788347507641437
$ws.Range("B8").Value = 0.7071695230269043
